$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.216.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6996"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.53"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07886"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.51"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08168"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.856.32"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.177"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7025"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.38"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.228.63"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.798"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007802"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.86"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.095.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.493"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.63"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.843"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1414"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.908"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.404"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.471"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.52%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05142"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.161"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7073"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9994"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.683"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01842"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.706"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.153.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9281"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.956"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4233"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5293"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.730"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.111"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.933"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.27%  "
